$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("ENUM(") -and $t.Contains("member") -and $t.Contains("admin")) {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $full = $d.Range($pStart, $pEnd)

        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2B8E4185" w14:textId="77777777" w:rsidR="00612C15" w:rsidRPr="00612C15" w:rsidRDefault="00612C15" w:rsidP="00612C15"><w:pPr><w:tabs><w:tab w:val="left" w:pos="7462"/></w:tabs></w:pPr><w:r w:rsidRPr="00612C15"><w:t>`role` ENUM(''</w:t></w:r><w:r><w:t>MENBER</w:t></w:r><w:r><w:t>'',''</w:t></w:r><w:r><w:t>ADMIN</w:t></w:r><w:r><w:t>'') NOT NULL DEFAULT ''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00612C15"><w:t>member</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00612C15"><w:t>'',</w:t></w:r></w:p>'

        $full.InsertXML($xml)
        Write-Output "Replaced paragraph"
        break
    }
}
Write-Output "Done"
